# Update scripts with new TPM values.
#
# The underlying analysis was re-run with new TPM-based expression data.
# This changed:
#   - the numeric statistics for the 4 remaining Sending-cluster rows
#     (each row now always targets the "Resolving-Mac" cluster, and the
#     "Target cluster" value for row 2 / row 4 changed accordingly)
#   - the (ECs, Resolving-Mac)/(FAPs, Resolving-Mac) pairing rows that
#     used to live in rows 6-9 are no longer present, so rows 6-9 are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Icosl -> Cd28 -> Resolving-Mac -----------------------
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 1.089622333333333
$ws.Range("H2").Value = 3.268867
$ws.Range("I2").Value = 0.09062709179941439
$ws.Range("J2").Value = 0.09062709179941439
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.273410666666667
$ws.Range("N2").Value = 15.820232
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 5.746026035238224
$ws.Range("R2").Value = 51.714234317144
$ws.Range("S2").Value = 0.09062709179941439
$ws.Range("T2").Value = 0.09062709179941439

# --- Row 3: FAPs -> Icosl -> Cd28 -> Resolving-Mac ----------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 1.876575666666667
$ws.Range("H3").Value = 5.629727
$ws.Range("I3").Value = 0.1560803133424032
$ws.Range("J3").Value = 0.1560803133424032
$ws.Range("M3").Value = 5.273410666666667
$ws.Range("N3").Value = 15.820232
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 9.895954137407111
$ws.Range("R3").Value = 89.06358723666401
$ws.Range("S3").Value = 0.1560803133424032
$ws.Range("T3").Value = 0.1560803133424032

# --- Row 4: MuSCs -> Icosl -> Cd28 -> Resolving-Mac ---------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 0.6928603333333333
$ws.Range("H4").Value = 2.078581
$ws.Range("I4").Value = 0.05762723019918477
$ws.Range("J4").Value = 0.05762723019918477
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.273410666666667
$ws.Range("N4").Value = 15.820232
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 3.653737072310222
$ws.Range("R4").Value = 32.883633650792
$ws.Range("S4").Value = 0.05762723019918477
$ws.Range("T4").Value = 0.05762723019918477

# --- Row 5: Resolving-Mac -> Icosl -> Cd28 -> Resolving-Mac -------------
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 8.364082999999999
$ws.Range("H5").Value = 25.092249
$ws.Range("I5").Value = 0.6956653646589975
$ws.Range("J5").Value = 0.6956653646589976
$ws.Range("M5").Value = 5.273410666666667
$ws.Range("N5").Value = 15.820232
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 44.10724450908533
$ws.Range("R5").Value = 396.965200581768
$ws.Range("S5").Value = 0.6956653646589975
$ws.Range("T5").Value = 0.6956653646589976

# --- Remove the now-obsolete rows 6-9 -----------------------------------
$ws.Range("A6:T9").EntireRow.Delete()
